$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.180165333333334
$ws.Cells.Item(2, 8).Value = 6.540496
$ws.Cells.Item(2, 9).Value = 0.01970539991828544
$ws.Cells.Item(2, 10).Value = 0.01970539991828544
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.4702473333333333
$ws.Cells.Item(2, 14).Value = 1.410742
$ws.Cells.Item(2, 15).Value = 0.009034922268422819
$ws.Cells.Item(2, 16).Value = 0.009034922268422819
$ws.Cells.Item(2, 17).Value = 1.025216934225778
$ws.Cells.Item(2, 18).Value = 9.226952408032
$ws.Cells.Item(2, 19).Value = 0.0001780367565298943
$ws.Cells.Item(2, 20).Value = 0.0001780367565298943
$ws.Cells.Item(3, 7).Value = 2.180165333333334
$ws.Cells.Item(3, 8).Value = 6.540496
$ws.Cells.Item(3, 9).Value = 0.01970539991828544
$ws.Cells.Item(3, 10).Value = 0.01970539991828544
$ws.Cells.Item(3, 14).Value = 0.9584440000000001
$ws.Cells.Item(3, 15).Value = 0.006138235792679485
$ws.Cells.Item(3, 16).Value = 0.006138235792679485
$ws.Cells.Item(3, 17).Value = 0.6965221275804445
$ws.Cells.Item(3, 18).Value = 6.268699148224001
$ws.Cells.Item(3, 19).Value = 0.0001209563910874831
$ws.Cells.Item(3, 20).Value = 0.0001209563910874831
$ws.Cells.Item(4, 7).Value = 2.180165333333334
$ws.Cells.Item(4, 8).Value = 6.540496
$ws.Cells.Item(4, 9).Value = 0.01970539991828544
$ws.Cells.Item(4, 10).Value = 0.01970539991828544
$ws.Cells.Item(4, 13).Value = 1.047307
$ws.Cells.Item(4, 14).Value = 3.141921
$ws.Cells.Item(4, 15).Value = 0.02012204358311108
$ws.Cells.Item(4, 16).Value = 0.02012204358311108
$ws.Cells.Item(4, 17).Value = 2.283302414757333
$ws.Cells.Item(4, 18).Value = 20.549721732816
$ws.Cells.Item(4, 19).Value = 0.0003965129159783731
$ws.Cells.Item(4, 20).Value = 0.0003965129159783731
$ws.Cells.Item(5, 7).Value = 2.180165333333334
$ws.Cells.Item(5, 8).Value = 6.540496
$ws.Cells.Item(5, 9).Value = 0.01970539991828544
$ws.Cells.Item(5, 10).Value = 0.01970539991828544
$ws.Cells.Item(5, 13).Value = 50.21070966666667
$ws.Cells.Item(5, 14).Value = 150.632129
$ws.Cells.Item(5, 15).Value = 0.9647047983557866
$ws.Cells.Item(5, 16).Value = 0.9647047983557866
$ws.Cells.Item(5, 17).Value = 109.4676485773316
$ws.Cells.Item(5, 18).Value = 985.2088371959842
$ws.Cells.Item(5, 19).Value = 0.01900989385468968
$ws.Cells.Item(5, 20).Value = 0.01900989385468968
$ws.Cells.Item(6, 9).Value = 0.733713204346044
$ws.Cells.Item(6, 10).Value = 0.7337132043460441
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.4702473333333333
$ws.Cells.Item(6, 14).Value = 1.410742
$ws.Cells.Item(6, 15).Value = 0.009034922268422819
$ws.Cells.Item(6, 16).Value = 0.009034922268422819
$ws.Cells.Item(6, 17).Value = 38.17304926973911
$ws.Cells.Item(6, 18).Value = 343.557443427652
$ws.Cells.Item(6, 19).Value = 0.006629041768581936
$ws.Cells.Item(6, 20).Value = 0.006629041768581937
$ws.Cells.Item(7, 9).Value = 0.733713204346044
$ws.Cells.Item(7, 10).Value = 0.7337132043460441
$ws.Cells.Item(7, 14).Value = 0.9584440000000001
$ws.Cells.Item(7, 15).Value = 0.006138235792679485
$ws.Cells.Item(7, 16).Value = 0.006138235792679485
$ws.Cells.Item(7, 17).Value = 25.93438774367378
$ws.Cells.Item(7, 19).Value = 0.004503704652478444
$ws.Cells.Item(7, 20).Value = 0.004503704652478445
$ws.Cells.Item(8, 9).Value = 0.733713204346044
$ws.Cells.Item(8, 10).Value = 0.7337132043460441
$ws.Cells.Item(8, 13).Value = 1.047307
$ws.Cells.Item(8, 14).Value = 3.141921
$ws.Cells.Item(8, 15).Value = 0.02012204358311108
$ws.Cells.Item(8, 16).Value = 0.02012204358311108
$ws.Cells.Item(8, 17).Value = 85.01675369034733
$ws.Cells.Item(8, 18).Value = 765.150783213126
$ws.Cells.Item(8, 19).Value = 0.01476380907535519
$ws.Cells.Item(8, 20).Value = 0.01476380907535519
$ws.Cells.Item(9, 9).Value = 0.733713204346044
$ws.Cells.Item(9, 10).Value = 0.7337132043460441
$ws.Cells.Item(9, 13).Value = 50.21070966666667
$ws.Cells.Item(9, 14).Value = 150.632129
$ws.Cells.Item(9, 15).Value = 0.9647047983557866
$ws.Cells.Item(9, 16).Value = 0.9647047983557866
$ws.Cells.Item(9, 17).Value = 4075.931447367909
$ws.Cells.Item(9, 18).Value = 36683.38302631118
$ws.Cells.Item(9, 19).Value = 0.7078166488496285
$ws.Cells.Item(9, 20).Value = 0.7078166488496286
$ws.Cells.Item(10, 7).Value = 25.672264
$ws.Cells.Item(10, 8).Value = 77.016792
$ws.Cells.Item(10, 9).Value = 0.2320384702908474
$ws.Cells.Item(10, 10).Value = 0.2320384702908474
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.4702473333333333
$ws.Cells.Item(10, 14).Value = 1.410742
$ws.Cells.Item(10, 15).Value = 0.009034922268422819
$ws.Cells.Item(10, 16).Value = 0.009034922268422819
$ws.Cells.Item(10, 17).Value = 12.07231368662933
$ws.Cells.Item(10, 18).Value = 108.650823179664
$ws.Cells.Item(10, 19).Value = 0.002096449542361544
$ws.Cells.Item(10, 20).Value = 0.002096449542361544
$ws.Cells.Item(11, 7).Value = 25.672264
$ws.Cells.Item(11, 8).Value = 77.016792
$ws.Cells.Item(11, 9).Value = 0.2320384702908474
$ws.Cells.Item(11, 10).Value = 0.2320384702908474
$ws.Cells.Item(11, 14).Value = 0.9584440000000001
$ws.Cells.Item(11, 15).Value = 0.006138235792679485
$ws.Cells.Item(11, 16).Value = 0.006138235792679485
$ws.Cells.Item(11, 17).Value = 8.201809132405334
$ws.Cells.Item(11, 18).Value = 73.816282191648
$ws.Cells.Item(11, 19).Value = 0.001424306843617875
$ws.Cells.Item(11, 20).Value = 0.001424306843617875
$ws.Cells.Item(12, 7).Value = 25.672264
$ws.Cells.Item(12, 8).Value = 77.016792
$ws.Cells.Item(12, 9).Value = 0.2320384702908474
$ws.Cells.Item(12, 10).Value = 0.2320384702908474
$ws.Cells.Item(12, 13).Value = 1.047307
$ws.Cells.Item(12, 14).Value = 3.141921
$ws.Cells.Item(12, 15).Value = 0.02012204358311108
$ws.Cells.Item(12, 16).Value = 0.02012204358311108
$ws.Cells.Item(12, 17).Value = 26.886741793048
$ws.Cells.Item(12, 18).Value = 241.980676137432
$ws.Cells.Item(12, 19).Value = 0.004669088212150858
$ws.Cells.Item(12, 20).Value = 0.004669088212150858
$ws.Cells.Item(13, 7).Value = 25.672264
$ws.Cells.Item(13, 8).Value = 77.016792
$ws.Cells.Item(13, 9).Value = 0.2320384702908474
$ws.Cells.Item(13, 10).Value = 0.2320384702908474
$ws.Cells.Item(13, 13).Value = 50.21070966666667
$ws.Cells.Item(13, 14).Value = 150.632129
$ws.Cells.Item(13, 15).Value = 0.9647047983557866
$ws.Cells.Item(13, 16).Value = 0.9647047983557866
$ws.Cells.Item(13, 17).Value = 1289.022594190019
$ws.Cells.Item(13, 18).Value = 11601.20334771017
$ws.Cells.Item(13, 19).Value = 0.2238486256927171
$ws.Cells.Item(13, 20).Value = 0.2238486256927171
$ws.Cells.Item(14, 7).Value = 1.608999666666667
$ws.Cells.Item(14, 8).Value = 4.826999
$ws.Cells.Item(14, 9).Value = 0.01454292544482312
$ws.Cells.Item(14, 10).Value = 0.01454292544482312
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.4702473333333333
$ws.Cells.Item(14, 14).Value = 1.410742
$ws.Cells.Item(14, 15).Value = 0.009034922268422819
$ws.Cells.Item(14, 16).Value = 0.009034922268422819
$ws.Cells.Item(14, 17).Value = 0.7566278025842221
$ws.Cells.Item(14, 18).Value = 6.809650223257999
$ws.Cells.Item(14, 19).Value = 0.0001313942009494453
$ws.Cells.Item(14, 20).Value = 0.0001313942009494453
$ws.Cells.Item(15, 7).Value = 1.608999666666667
$ws.Cells.Item(15, 8).Value = 4.826999
$ws.Cells.Item(15, 9).Value = 0.01454292544482312
$ws.Cells.Item(15, 10).Value = 0.01454292544482312
$ws.Cells.Item(15, 14).Value = 0.9584440000000001
$ws.Cells.Item(15, 15).Value = 0.006138235792679485
$ws.Cells.Item(15, 16).Value = 0.006138235792679485
$ws.Cells.Item(15, 17).Value = 0.5140453588395556
$ws.Cells.Item(15, 18).Value = 4.626408229556001
$ws.Cells.Item(15, 19).Value = 0.00008926790549568253
$ws.Cells.Item(15, 20).Value = 0.00008926790549568253
$ws.Cells.Item(16, 7).Value = 1.608999666666667
$ws.Cells.Item(16, 8).Value = 4.826999
$ws.Cells.Item(16, 9).Value = 0.01454292544482312
$ws.Cells.Item(16, 10).Value = 0.01454292544482312
$ws.Cells.Item(16, 13).Value = 1.047307
$ws.Cells.Item(16, 14).Value = 3.141921
$ws.Cells.Item(16, 15).Value = 0.02012204358311108
$ws.Cells.Item(16, 16).Value = 0.02012204358311108
$ws.Cells.Item(16, 17).Value = 1.685116613897667
$ws.Cells.Item(16, 18).Value = 15.166049525079
$ws.Cells.Item(16, 19).Value = 0.000292633379626666
$ws.Cells.Item(16, 20).Value = 0.000292633379626666
$ws.Cells.Item(17, 7).Value = 1.608999666666667
$ws.Cells.Item(17, 8).Value = 4.826999
$ws.Cells.Item(17, 9).Value = 0.01454292544482312
$ws.Cells.Item(17, 10).Value = 0.01454292544482312
$ws.Cells.Item(17, 13).Value = 50.21070966666667
$ws.Cells.Item(17, 14).Value = 150.632129
$ws.Cells.Item(17, 15).Value = 0.9647047983557866
$ws.Cells.Item(17, 16).Value = 0.9647047983557866
$ws.Cells.Item(17, 17).Value = 80.78901511676345
$ws.Cells.Item(17, 18).Value = 727.1011360508711
$ws.Cells.Item(17, 19).Value = 0.01402962995875133
$ws.Cells.Item(17, 20).Value = 0.01402962995875133
